$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates pulled from the "Updated symbol list" refresh.
# Columns B/C (coin name + link) are plain text already, so a direct
# .Value assignment keeps them as text.
# Columns D/E hold numeric-looking strings (prices / percentages) that
# must stay text (e.g. "277.80", "0.06367", "1.14%") rather than being
# auto-coerced to numbers, so those cells are pre-formatted as text via
# NumberFormat "@" before the value is written.
$textColumns = @("D", "E")

$updates = @(
    @{ Ref = "D2"; Value = "277.80" }
    @{ Ref = "E2"; Value = "0.89%" }
    @{ Ref = "E3"; Value = "0.01%" }
    @{ Ref = "E4"; Value = "2.07%" }
    @{ Ref = "D5"; Value = "0.06367" }
    @{ Ref = "E5"; Value = "0.83%" }
    @{ Ref = "D6"; Value = "7.029" }
    @{ Ref = "E6"; Value = "1.48%" }
    @{ Ref = "D7"; Value = "1.323" }
    @{ Ref = "E7"; Value = "1.21%" }
    @{ Ref = "D8"; Value = "0.8978" }
    @{ Ref = "E8"; Value = "2.23%" }
    @{ Ref = "D9"; Value = "0.1532" }
    @{ Ref = "E9"; Value = "0.70%" }
    @{ Ref = "B10"; Value = "LiechtensteinCryptoassetsExchange" }
    @{ Ref = "C10"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" }
    @{ Ref = "D10"; Value = "0.05208" }
    @{ Ref = "E10"; Value = "4.48%" }
    @{ Ref = "B11"; Value = "MandalaExchangeToken" }
    @{ Ref = "C11"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" }
    @{ Ref = "D11"; Value = "0.07492" }
    @{ Ref = "E11"; Value = "-0.21%" }
    @{ Ref = "B12"; Value = "BitrueCoin" }
    @{ Ref = "C12"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Ref = "D12"; Value = "0.02905" }
    @{ Ref = "E12"; Value = "0.31%" }
    @{ Ref = "B13"; Value = "BitMartToken" }
    @{ Ref = "C13"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" }
    @{ Ref = "D13"; Value = "0.08961" }
    @{ Ref = "E13"; Value = "-1.03%" }
    @{ Ref = "B14"; Value = "BitForexToken" }
    @{ Ref = "C14"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" }
    @{ Ref = "D14"; Value = "0.001568" }
    @{ Ref = "E14"; Value = "0.00%" }
    @{ Ref = "B15"; Value = "One" }
    @{ Ref = "C15"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" }
    @{ Ref = "D15"; Value = "0.0006390" }
    @{ Ref = "E15"; Value = "0.78%" }
    @{ Ref = "B16"; Value = "TigerCash" }
    @{ Ref = "C16"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Ref = "D16"; Value = "0.006048" }
    @{ Ref = "E16"; Value = "1.88%" }
    @{ Ref = "B17"; Value = "LEO" }
    @{ Ref = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Ref = "D17"; Value = "3.480" }
    @{ Ref = "E17"; Value = "0.80%" }
    @{ Ref = "B18"; Value = "GateToken" }
    @{ Ref = "C18"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" }
    @{ Ref = "D18"; Value = "3.303" }
    @{ Ref = "E18"; Value = "0.10%" }
    @{ Ref = "B19"; Value = "BTSEToken" }
    @{ Ref = "C19"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" }
    @{ Ref = "D19"; Value = "2.247" }
    @{ Ref = "E19"; Value = "-1.09%" }
    @{ Ref = "B20"; Value = "BitpandaEcosystemToken" }
    @{ Ref = "C20"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best" }
    @{ Ref = "D20"; Value = "0.3093" }
    @{ Ref = "E20"; Value = "-0.82%" }
    @{ Ref = "D21"; Value = "0.1332" }
    @{ Ref = "E21"; Value = "0.49%" }
    @{ Ref = "D22"; Value = "3.909" }
    @{ Ref = "E22"; Value = "0.10%" }
    @{ Ref = "E23"; Value = "11.31%" }
    @{ Ref = "D24"; Value = "0.04393" }
    @{ Ref = "E24"; Value = "-0.11%" }
    @{ Ref = "D25"; Value = "0.001171" }
    @{ Ref = "E25"; Value = "0.02%" }
    @{ Ref = "E26"; Value = "1.48%" }
    @{ Ref = "D28"; Value = "0.0001181" }
    @{ Ref = "D29"; Value = "0.0001649" }
    @{ Ref = "E29"; Value = "-14.88%" }
    @{ Ref = "D40"; Value = "0.04071" }
    @{ Ref = "E40"; Value = "-1.10%" }
    @{ Ref = "D41"; Value = "0.006794" }
    @{ Ref = "E41"; Value = "-1.15%" }
    @{ Ref = "D42"; Value = "0.1410" }
    @{ Ref = "E42"; Value = "20.11%" }
    @{ Ref = "E43"; Value = "-2.18%" }
    @{ Ref = "D44"; Value = "0.01166" }
    @{ Ref = "E44"; Value = "1.44%" }
    @{ Ref = "D45"; Value = "0.00005350" }
    @{ Ref = "E45"; Value = "1.97%" }
    @{ Ref = "D46"; Value = "1.561" }
    @{ Ref = "E46"; Value = "4.79%" }
    @{ Ref = "D47"; Value = "0.01851" }
    @{ Ref = "E47"; Value = "-7.45%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $colLetter = $u.Ref -replace '[0-9]+$', ''
    if ($textColumns -contains $colLetter) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
